$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Dropdown A
$ws.Range("A2").Value = "fa694862-f183-4200-8896-4caf5cec47d3"
$ws.Range("B2").Value = "Dropdown A"
$ws.Range("C2").Value = "A high-qauality dropdown"
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 23

# Row 3 - Dropdown B
$ws.Range("A3").Value = "7a95839e-7075-40c8-9c46-a5990084fb46"
$ws.Range("B3").Value = "Dropdown B"
$ws.Range("C3").Value = "A high-qauality dropdown"
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 25

# Row 4 - Widget D
$ws.Range("A4").Value = "ba6f6cb7-a21f-4898-a1df-2731b4239fcb"
$ws.Range("B4").Value = "Widget D"
$ws.Range("C4").Value = "A good-qauality widget"
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 24

# Row 5 - Widget C
$ws.Range("A5").Value = "e0f79033-b1c3-437d-9b29-95057c11baa2"
$ws.Range("B5").Value = "Widget C"
$ws.Range("C5").Value = "A premium widget"
$ws.Range("D5").Value = 35
$ws.Range("E5").Value = 30

# Row 6 - new - Widget E
$ws.Range("A6").Value = "9a3c15bf-2cd8-4c06-9c43-a1aaf680c68a"
$ws.Range("B6").Value = "Widget E"
$ws.Range("C6").Value = "A heavily used widget"
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 10

$ws.Rows("6:6").Select() | Out-Null
